# "Remove Data Row Activities"
#
# This workbook is used as sample data for a UiPath "DataTableActivity" demo.
# The commit removes a data row from the "Sheet4" worksheet (simulating the
# result of a "Remove Data Row" activity that dropped the 1014/Salman row and
# left two trailing duplicate rows of the last remaining entry), and clears
# out a stray duplicate pair of columns (F:G) that had been written onto the
# "Sheet3" worksheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet3: remove the duplicated "Student Name"/"Roll No." columns (F:G)
# that mirrored columns A and C. Clearing their contents shrinks the
# worksheet's used range back down to A1:C15.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("F1:G15").ClearContents()
$ws3.Activate()
$ws3.Range("E1:H15").Select()

# ---------------------------------------------------------------------
# Sheet4: remove the data row for Roll No. 1014 (Salman / Business).
# That deletion shifts the remaining rows up, so what used to be row 15
# (Roll No. 1017, Yi Wong / Micro Biology) is copied back in twice to
# restore row count/dimension to 15 rows, matching the activity's output.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Rows(12).Delete()

$lastRow = $ws4.Range("A13:G13")
$lastRow.Copy($ws4.Range("A14:G14"))
$lastRow.Copy($ws4.Range("A15:G15"))

$ws4.Activate()
$ws4.Range("A1:G15").Select()

# ---------------------------------------------------------------------
# Sheet2: no data on this sheet; only the saved cursor/selection moved.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()
$ws2.Range("I19").Select()

# Re-activate Sheet4 last, matching the workbook's active tab.
$ws4.Activate()
